$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "64.291.20"
$ws.Range("E2").Value2 = "  -0.31%  "
$ws.Range("D3").Value2 = "3.136.96"
$ws.Range("E3").Value2 = "  -1.36%  "
$ws.Range("E4").Value2 = "  -0.04%  "
$ws.Range("D5").Value2 = "'571.84"
$ws.Range("E5").Value2 = "  +0.00%  "
$ws.Range("D6").Value2 = "'163.71"
$ws.Range("E6").Value2 = "  -3.20%  "
$ws.Range("E8").Value2 = "  -5.33%  "
$ws.Range("D9").Value2 = "3.149.42"
$ws.Range("E9").Value2 = "  -1.28%  "
$ws.Range("E10").Value2 = "  -3.37%  "
$ws.Range("E11").Value2 = "  -3.08%  "
$ws.Range("E12").Value2 = "  +0.04%  "
$ws.Range("D13").Value2 = "3.685.38"
$ws.Range("E13").Value2 = "  -1.46%  "
$ws.Range("E14").Value2 = "  -1.42%  "
$ws.Range("D15").Value2 = "64.343.04"
$ws.Range("E15").Value2 = "  -0.30%  "
$ws.Range("D16").Value2 = "'25.12"
$ws.Range("E16").Value2 = "  -1.15%  "
$ws.Range("D17").Value2 = "3.151.21"
$ws.Range("E17").Value2 = "  -1.31%  "
$ws.Range("E18").Value2 = "  -3.10%  "
$ws.Range("D19").Value2 = "'401.88"
$ws.Range("E19").Value2 = "  -3.70%  "
$ws.Range("E20").Value2 = "  -1.84%  "
$ws.Range("D21").Value2 = "'12.51"
$ws.Range("E21").Value2 = "  -3.34%  "
$ws.Range("E22").Value2 = "  -0.56%  "
$ws.Range("D23").Value2 = "'0.999"
$ws.Range("E23").Value2 = "  -0.05%  "
$ws.Range("D24").Value2 = "'68.70"
$ws.Range("E24").Value2 = "  -2.66%  "
$ws.Range("E25").Value2 = "  -0.88%  "
$ws.Range("E26").Value2 = "  -4.16%  "
$ws.Range("E27").Value2 = "  -4.49%  "
$ws.Range("D28").Value2 = "'8.77"
$ws.Range("E28").Value2 = "  -0.14%  "
$ws.Range("D29").Value2 = "'0.997"
$ws.Range("E29").Value2 = "  -0.05%  "
$ws.Range("E31").Value2 = "  -1.61%  "
$ws.Range("E32").Value2 = "  -2.68%  "
$ws.Range("D33").Value2 = "'161.24"
$ws.Range("E33").Value2 = "  +2.21%  "
$ws.Range("B34").Value2 = "NEARProtocol"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value2 = "'4.85"
$ws.Range("E34").Value2 = "  -4.90%  "
$ws.Range("B35").Value2 = "Aptos"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value2 = "'6.27"
$ws.Range("E35").Value2 = "  -1.03%  "
$ws.Range("E36").Value2 = "  -2.45%  "
$ws.Range("E37").Value2 = "  -1.70%  "
$ws.Range("E38").Value2 = "  -1.86%  "
$ws.Range("D39").Value2 = "2.636.52"
$ws.Range("E39").Value2 = "  -3.69%  "
$ws.Range("D40").Value2 = "'23.71"
$ws.Range("E40").Value2 = "  -3.11%  "
$ws.Range("D41").Value2 = "'4.07"
$ws.Range("E41").Value2 = "  -3.33%  "
$ws.Range("D42").Value2 = "'38.38"
$ws.Range("E42").Value2 = "  -2.20%  "
$ws.Range("D43").Value2 = "'0.689"
$ws.Range("E43").Value2 = "  -3.98%  "
$ws.Range("D44").Value2 = "'0.0613"
$ws.Range("E44").Value2 = "  -1.65%  "
$ws.Range("D45").Value2 = "'5.43"
$ws.Range("E45").Value2 = "  -4.57%  "
$ws.Range("E46").Value2 = "  -3.90%  "
$ws.Range("E47").Value2 = "  -2.79%  "
$ws.Range("D48").Value2 = "'286.22"
$ws.Range("E48").Value2 = "  -2.76%  "
$ws.Range("D49").Value2 = "'0.996"
$ws.Range("E49").Value2 = "  -0.12%  "
$ws.Range("D50").Value2 = "'0.0976"
$ws.Range("E50").Value2 = "  -1.35%  "
$ws.Range("D51").Value2 = "'10.49"
$ws.Range("E51").Value2 = "  +0.52%  "
